$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("dataset_income_2018-2023")

# The "percent of total" formula rows (142, 144, ... 160) on the
# dataset_income_2018-2023 sheet each multiply the year total in row 140
# by a percentage figure taken from a neighboring data row, divided by 100.
# They previously pointed one row BELOW themselves (row+1); they now need
# to point one row ABOVE themselves (row-1). Re-entering the formula on
# each row (C:BB) lets Excel's relative-reference fill logic recompute the
# per-column reference and the cached values for every column in one shot.
for ($row = 142; $row -le 160; $row += 2) {
    $targetRow = $row - 1
    $ws.Range("C$row").Formula = "=C`$140*C$targetRow/100"
    $ws.Range("D$row`:BB$row").Formula = "=D`$140*D$targetRow/100"
}

# Restore the reported scroll position / active selection for the sheet.
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 157
$win.ScrollColumn = 1
$ws.Range("C167").Select()
